$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A3 from "M2x20" to "M3x20"
$ws.Range("A3").Value = "M3x20"

# Update the active selection to K17 (matches post-edit cursor position in diff)
$ws.Range("K17").Select()
